{"js": "// Lattice-multiplication worksheet refresh: swap in the new set of\n// multiplication problems (and their partial-product lattice digits)\n// cell-by-cell, keeping each cell's existing run formatting (sz=32).\n//\n// The table is a fixed 5-row x 3-column grid; only the text inside each\n// cell changes (the \"x\" problem line, the two factor digits, the divider,\n// and the two leading lattice digits). \"\\u000B\" (vertical tab) is Word's\n// manual-line-break character, which insertText() turns into a <w:br/>\n// between <w:t> runs, matching the original markup shape.\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\nconst table = tables.items[0];\n\nconst cellValues = [\n  { row: 0, col: 0, text: \"89 x 73\\u000B  7    3\\u000B  ----\\u000B8|    |\\u000B9|    |\" },\n  { row: 0, col: 1, text: \"45 x 11\\u000B  1    1\\u000B  ----\\u000B4|    |\\u000B5|    |\" },\n  { row: 0, col: 2, text: \"59 x 69\\u000B  6    9\\u000B  ----\\u000B5|    |\\u000B9|    |\" },\n  { row: 1, col: 0, text: \"89 x 15\\u000B  1    5\\u000B  ----\\u000B8|    |\\u000B9|    |\" },\n  { row: 1, col: 1, text: \"99 x 79\\u000B  7    9\\u000B  ----\\u000B9|    |\\u000B9|    |\" },\n  { row: 1, col: 2, text: \"98 x 20\\u000B  2    0\\u000B  ----\\u000B9|    |\\u000B8|    |\" },\n  { row: 2, col: 0, text: \"97 x 39\\u000B  3    9\\u000B  ----\\u000B9|    |\\u000B7|    |\" },\n  { row: 2, col: 1, text: \"26 x 85\\u000B  8    5\\u000B  ----\\u000B2|    |\\u000B6|    |\" },\n  { row: 2, col: 2, text: \"50 x 86\\u000B  8    6\\u000B  ----\\u000B5|    |\\u000B0|    |\" },\n  { row: 3, col: 0, text: \"94 x 92\\u000B  9    2\\u000B  ----\\u000B9|    |\\u000B4|    |\" },\n  { row: 3, col: 1, text: \"77 x 44\\u000B  4    4\\u000B  ----\\u000B7|    |\\u000B7|    |\" },\n  { row: 3, col: 2, text: \"70 x 64\\u000B  6    4\\u000B  ----\\u000B7|    |\\u000B0|    |\" },\n  { row: 4, col: 0, text: \"83 x 60\\u000B  6    0\\u000B  ----\\u000B8|    |\\u000B3|    |\" },\n  { row: 4, col: 1, text: \"35 x 30\\u000B  3    0\\u000B  ----\\u000B3|    |\\u000B5|    |\" },\n  { row: 4, col: 2, text: \"10 x 85\\u000B  8    5\\u000B  ----\\u000B1|    |\\u000B0|    |\" },\n];\n\nfor (const { row, col, text } of cellValues) {\n  const cell = table.getCell(row, col);\n  const rng = cell.body.getRange();\n  rng.insertText(text, Word.InsertLocation.replace);\n}\nawait context.sync();\n", "ps1": "# Lattice-multiplication worksheet refresh: swap in the new set of\n# multiplication problems (and their partial-product lattice digits)\n# cell-by-cell, keeping each cell's existing run formatting (sz=32).\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n$nl = [char]11   # vertical-tab = Word's \"manual line break\" (<w:br/>) in Range.Text\n\n$t.Cell(1, 1).Range.Text = \"89 x 73\" + $nl + \"  7    3\" + $nl + \"  ----\" + $nl + \"8|    |\" + $nl + \"9|    |\"\n$t.Cell(1, 2).Range.Text = \"45 x 11\" + $nl + \"  1    1\" + $nl + \"  ----\" + $nl + \"4|    |\" + $nl + \"5|    |\"\n$t.Cell(1, 3).Range.Text = \"59 x 69\" + $nl + \"  6    9\" + $nl + \"  ----\" + $nl + \"5|    |\" + $nl + \"9|    |\"\n\n$t.Cell(2, 1).Range.Text = \"89 x 15\" + $nl + \"  1    5\" + $nl + \"  ----\" + $nl + \"8|    |\" + $nl + \"9|    |\"\n$t.Cell(2, 2).Range.Text = \"99 x 79\" + $nl + \"  7    9\" + $nl + \"  ----\" + $nl + \"9|    |\" + $nl + \"9|    |\"\n$t.Cell(2, 3).Range.Text = \"98 x 20\" + $nl + \"  2    0\" + $nl + \"  ----\" + $nl + \"9|    |\" + $nl + \"8|    |\"\n\n$t.Cell(3, 1).Range.Text = \"97 x 39\" + $nl + \"  3    9\" + $nl + \"  ----\" + $nl + \"9|    |\" + $nl + \"7|    |\"\n$t.Cell(3, 2).Range.Text = \"26 x 85\" + $nl + \"  8    5\" + $nl + \"  ----\" + $nl + \"2|    |\" + $nl + \"6|    |\"\n$t.Cell(3, 3).Range.Text = \"50 x 86\" + $nl + \"  8    6\" + $nl + \"  ----\" + $nl + \"5|    |\" + $nl + \"0|    |\"\n\n$t.Cell(4, 1).Range.Text = \"94 x 92\" + $nl + \"  9    2\" + $nl + \"  ----\" + $nl + \"9|    |\" + $nl + \"4|    |\"\n$t.Cell(4, 2).Range.Text = \"77 x 44\" + $nl + \"  4    4\" + $nl + \"  ----\" + $nl + \"7|    |\" + $nl + \"7|    |\"\n$t.Cell(4, 3).Range.Text = \"70 x 64\" + $nl + \"  6    4\" + $nl + \"  ----\" + $nl + \"7|    |\" + $nl + \"0|    |\"\n\n$t.Cell(5, 1).Range.Text = \"83 x 60\" + $nl + \"  6    0\" + $nl + \"  ----\" + $nl + \"8|    |\" + $nl + \"3|    |\"\n$t.Cell(5, 2).Range.Text = \"35 x 30\" + $nl + \"  3    0\" + $nl + \"  ----\" + $nl + \"3|    |\" + $nl + \"5|    |\"\n$t.Cell(5, 3).Range.Text = \"10 x 85\" + $nl + \"  8    5\" + $nl + \"  ----\" + $nl + \"1|    |\" + $nl + \"0|    |\"\n"}
